# AL KABIR TOWER - Opening Schedule
# Door dimensions corrected to use architect's text annotations (kapi-yazi
# layer) instead of xscale-derived values:
#   Single Leaf Hinged Door      880mm  -> 900mm
#   Single Leaf Hinged Door      980mm  -> 1000mm
#   Double Leaf Hinged Entrance Door 1980mm -> 2000mm
# Ripple: Total Doors 422->413, Total Openings 667->658, plus every
# dependent area / count cell across the six sheets.

# ---------------------------------------------------------------------
# Helper: write a string value into a cell while keeping the cell's
# existing style (number format, font, fill, borders, alignment...)
# untouched. A plain `$Range.Value = "413"` assignment would make Excel
# auto-coerce an all-digit string into a real number (changing the
# stored type from text to numeric), which is not what we want here —
# these "422" / "667" style cells are text labels, not numbers. Routing
# the write through a Text-formatted scratch cell + PasteSpecial(values
# only) keeps the literal text type without touching the destination's
# style.
function Set-TextValue {
    param($Range, $Text)
    $ws = $Range.Worksheet
    $scratch = $ws.Range("ZZ9999")
    $scratch.NumberFormat = "@"
    $scratch.Value = $Text
    $scratch.Copy() | Out-Null
    $Range.PasteSpecial(-4163) | Out-Null
    $scratch.Clear() | Out-Null
}

$wb = $excel.ActiveWorkbook

# =======================================================================
# Sheet: Project Summary
# =======================================================================
$ws1 = $wb.Worksheets.Item("Project Summary")

# SCHEDULE STATISTICS block (text-typed "label: value" cells)
Set-TextValue $ws1.Range("C19") "413"          # Total Doors:
Set-TextValue $ws1.Range("C20") "658"          # Total Openings:
$ws1.Range("C22").Value = "855.1 m²"           # Door Area:
$ws1.Range("C23").Value = "1799.8 m²"          # Total Facade Area:

# ALL OPENINGS OVERVIEW table - door rows
$ws1.Range("D36").Value = 900
$ws1.Range("G36").Value = 353
$ws1.Range("H36").Value = 1.98
$ws1.Range("I36").Value = 698.9

$ws1.Range("D37").Value = 1000
$ws1.Range("H37").Value = 2.2
$ws1.Range("I37").Value = 107.8

$ws1.Range("D38").Value = 2000
$ws1.Range("H38").Value = 4.4
$ws1.Range("I38").Value = 48.4

$ws1.Range("G39").Value = 658
$ws1.Range("I39").Value = 1799.8

# =======================================================================
# Sheet: Door Schedule
# =======================================================================
$ws3 = $wb.Worksheets.Item("Door Schedule")

$ws3.Range("C6").Value = 900
$ws3.Range("F6").Value = 37
$ws3.Range("H6").Value = 353
$ws3.Range("I6").Value = 1.98
$ws3.Range("J6").Value = 698.9

$ws3.Range("C7").Value = 1000
$ws3.Range("I7").Value = 2.2
$ws3.Range("J7").Value = 107.8

$ws3.Range("C8").Value = 2000
$ws3.Range("I8").Value = 4.4
$ws3.Range("J8").Value = 48.4

$ws3.Range("F9").Value = 43
$ws3.Range("H9").Value = 413
$ws3.Range("J9").Value = 855.1

# =======================================================================
# Sheet: Per-Floor Breakdown
# =======================================================================
$ws4 = $wb.Worksheets.Item("Per-Floor Breakdown")

# Row 18 - Single Leaf Hinged Door 900mm (per-floor 38 -> 37, total 362 -> 353)
$ws4.Range("D18").Value = 900
foreach ($col in @("F","G","H","I","J","K","L","M","N")) {
    $ws4.Range($col + "18").Value = 37
}
$ws4.Range("P18").Value = 353

# Row 19 - Single Leaf Hinged Door 1000mm (only width changes)
$ws4.Range("D19").Value = 1000

# Row 20 - Double Leaf Hinged Entrance Door 2000mm (only width changes)
$ws4.Range("D20").Value = 2000

# Row 21 - DOOR SUBTOTAL (422 -> 413)
foreach ($col in @("F","G","H","I","J","K","L","M","N")) {
    $ws4.Range($col + "21").Value = 43
}
$ws4.Range("P21").Value = 413

# Row 22 - GRAND TOTAL (667 -> 658)
foreach ($col in @("F","G","H","I","J","K","L","M","N")) {
    $ws4.Range($col + "22").Value = 68
}
$ws4.Range("P22").Value = 658

# =======================================================================
# Sheet: Area Summary
# =======================================================================
$ws5 = $wb.Worksheets.Item("Area Summary")

$ws5.Range("A19").Value = "DOORS (3 types, 413 units)"

$ws5.Range("D20").Value = 900
$ws5.Range("F20").Value = 1.98
$ws5.Range("G20").Value = 353
$ws5.Range("H20").Value = 698.9

$ws5.Range("D21").Value = 1000
$ws5.Range("F21").Value = 2.2
$ws5.Range("H21").Value = 107.8

$ws5.Range("D22").Value = 2000
$ws5.Range("F22").Value = 4.4
$ws5.Range("H22").Value = 48.4

$ws5.Range("G23").Value = 413
$ws5.Range("H23").Value = 855.1

$ws5.Range("G24").Value = 658
$ws5.Range("H24").Value = 1799.7

# =======================================================================
# Sheet: Quantity Verification
# =======================================================================
$ws6 = $wb.Worksheets.Item("Quantity Verification")

$ws6.Range("B30").Value = "Single Leaf Hinged Door 900mm"
$ws6.Range("D30").Value = 37
$ws6.Range("H30").Value = 333

$ws6.Range("I31").Value = 353

$ws6.Range("B32").Value = "Single Leaf Hinged Door 1000mm"

$ws6.Range("B34").Value = "Double Leaf Hinged Entrance Door 2000mm"

$ws6.Range("I36").Value = 413
$ws6.Range("I37").Value = 658
